$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42: 4600
$ws.Range("H42").Value = 438.875
$ws.Range("I42").Value = 557.1667
$ws.Range("K42").Value = 1671.5001
$ws.Range("M42").Value = -1441.5001
# Row 64: 5506
$ws.Range("H64").Value = 22467.455
$ws.Range("I64").Value = 23714.2
$ws.Range("K64").Value = 23714.2
$ws.Range("M64").Value = -23466.2
# Row 67: 5506
$ws.Range("H67").Value = 22467.455
$ws.Range("I67").Value = 23714.2
$ws.Range("K67").Value = 23714.2
$ws.Range("M67").Value = -22856.2
# Row 100: 19906
$ws.Range("H100").Value = 38028.234
$ws.Range("I100").Value = 51616.168
$ws.Range("J100").Value = 5417.2
$ws.Range("K100").Value = 51616.168
$ws.Range("L100").Value = 5417.2
$ws.Range("M100").Value = -51075.168
$ws.Range("N100").Value = -6499.2
# Row 125: 36228
$ws.Range("H125").Value = 38783.875
$ws.Range("I125").Value = 75530.25
$ws.Range("K125").Value = 679772.25
$ws.Range("M125").Value = -677312.25
# Row 138: 44169
$ws.Range("H138").Value = 3937.3447
$ws.Range("I138").Value = 3610
$ws.Range("J138").Value = 4340.231
$ws.Range("K138").Value = 10830
$ws.Range("L138").Value = 13020.693
$ws.Range("M138").Value = -5690
$ws.Range("N138").Value = -23300.693

$ws = $wb.Worksheets.Item("ARM")
# Row 18: 2470
$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20644
# Row 32: 44147
$ws.Range("H32").Value = 25480.63
$ws.Range("I32").Value = 28778.426
$ws.Range("K32").Value = 28778.426
$ws.Range("M32").Value = -28491.426
# Row 61: 43999
$ws.Range("H61").Value = 4214.4146
$ws.Range("I61").Value = 1508.3103
$ws.Range("K61").Value = 1508.3103
$ws.Range("M61").Value = -1296.3103
# Row 130: 34732
$ws.Range("H130").Value = 57422.5
$ws.Range("J130").Value = 57422.5
$ws.Range("L130").Value = 57422.5
$ws.Range("N130").Value = -67462.5
# Row 132: 43997
$ws.Range("H132").Value = 897
$ws.Range("I132").Value = 796.5
$ws.Range("K132").Value = 2389.5
$ws.Range("M132").Value = 140.5
# Row 136: 43999
$ws.Range("H136").Value = 4214.4146
$ws.Range("I136").Value = 1508.3103
$ws.Range("K136").Value = 4524.9309
$ws.Range("M136").Value = -1974.9309

$ws = $wb.Worksheets.Item("BSM")
# Row 86: 12526
$ws.Range("H86").Value = 1736.6666
$ws.Range("I86").Value = 1537.8889
$ws.Range("J86").Value = 2333
$ws.Range("K86").Value = 1537.8889
$ws.Range("L86").Value = 2333
$ws.Range("M86").Value = -414.8888999999999
$ws.Range("N86").Value = -4579
# Row 89: 12526
$ws.Range("H89").Value = 1736.6666
$ws.Range("I89").Value = 1537.8889
$ws.Range("J89").Value = 2333
$ws.Range("K89").Value = 7689.4445
$ws.Range("L89").Value = 11665
$ws.Range("M89").Value = -2073.4445
$ws.Range("N89").Value = -22897
# Row 99: 19943
$ws.Range("H99").Value = 1996.4166
$ws.Range("I99").Value = 1789.2858
$ws.Range("J99").Value = 2286.4
$ws.Range("K99").Value = 1789.2858
$ws.Range("L99").Value = 2286.4
$ws.Range("M99").Value = -291.2858000000001
$ws.Range("N99").Value = -5282.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31: 44023
$ws.Range("H31").Value = 5560231
$ws.Range("I31").Value = 10002880
$ws.Range("J31").Value = 6919.875
$ws.Range("K31").Value = 10002880
$ws.Range("L31").Value = 6919.875
$ws.Range("M31").Value = -10002585
$ws.Range("N31").Value = -7509.875
# Row 34: 44023
$ws.Range("H34").Value = 5560231
$ws.Range("I34").Value = 10002880
$ws.Range("J34").Value = 6919.875
$ws.Range("K34").Value = 10002880
$ws.Range("L34").Value = 6919.875
$ws.Range("M34").Value = -10002678
$ws.Range("N34").Value = -7323.875
# Row 122: 36196
$ws.Range("H122").Value = 3404.7144
$ws.Range("J122").Value = 3415.6
$ws.Range("L122").Value = 10246.8
$ws.Range("N122").Value = -15146.8

$ws = $wb.Worksheets.Item("CUL")
# Row 107: 27838
$ws.Range("H107").Value = 2471.375
$ws.Range("I107").Value = 7837
$ws.Range("J107").Value = 682.8333
$ws.Range("K107").Value = 23511
$ws.Range("L107").Value = 2048.4999
$ws.Range("M107").Value = -21591
$ws.Range("N107").Value = -5888.4999
# Row 110: 27857
$ws.Range("H110").Value = 12669.223
# Row 131: 36060
$ws.Range("H131").Value = 1741.75
$ws.Range("J131").Value = 1852
$ws.Range("L131").Value = 5556
$ws.Range("N131").Value = -15636

$ws = $wb.Worksheets.Item("GSM")
# Row 132: 44008
$ws.Range("H132").Value = 2040.8507
$ws.Range("I132").Value = 2114.6064
$ws.Range("J132").Value = 1291
$ws.Range("K132").Value = 6343.8192
$ws.Range("L132").Value = 3873
$ws.Range("M132").Value = -3813.8192
$ws.Range("N132").Value = -8933

$ws = $wb.Worksheets.Item("LTW")
# Row 17: 3017
$ws.Range("H17").Value = 79999
$ws.Range("J17").Value = 79999
$ws.Range("L17").Value = 79999
$ws.Range("N17").Value = -80339
# Row 22: 5277
$ws.Range("H22").Value = 1047.1177
$ws.Range("I22").Value = 816.75
$ws.Range("K22").Value = 816.75
$ws.Range("M22").Value = -521.75
# Row 27: 5277
$ws.Range("H27").Value = 1047.1177
$ws.Range("I27").Value = 816.75
$ws.Range("K27").Value = 816.75
$ws.Range("M27").Value = -709.75
# Row 36: 34261
$ws.Range("H36").Value = 75000
$ws.Range("J36").Value = 75000
$ws.Range("L36").Value = 75000
$ws.Range("N36").Value = -76124
# Row 61: 27740
$ws.Range("H61").Value = 877.8889
$ws.Range("I61").Value = 667.46155
$ws.Range("K61").Value = 667.46155
$ws.Range("M61").Value = -465.46155
# Row 100: 19995
$ws.Range("H100").Value = 2724.25
$ws.Range("I100").Value = 2501.0833
$ws.Range("J100").Value = 3393.75
$ws.Range("K100").Value = 2501.0833
$ws.Range("L100").Value = 3393.75
$ws.Range("M100").Value = -1960.0833
$ws.Range("N100").Value = -4475.75
# Row 113: 27740
$ws.Range("H113").Value = 877.8889
$ws.Range("I113").Value = 667.46155
$ws.Range("K113").Value = 667.46155
$ws.Range("M113").Value = 1502.53845
# Row 122: 36247
$ws.Range("H122").Value = 6165.278
$ws.Range("I122").Value = 5281.6665
$ws.Range("J122").Value = 7048.8887
$ws.Range("K122").Value = 15844.9995
$ws.Range("L122").Value = 21146.6661
$ws.Range("M122").Value = -13394.9995
$ws.Range("N122").Value = -26046.6661
# Row 132: 44058
$ws.Range("H132").Value = 2414.525
$ws.Range("I132").Value = 1902.6471
$ws.Range("K132").Value = 5707.9413
$ws.Range("M132").Value = -3177.9413
# Row 139: 43310
$ws.Range("H139").Value = 79994
$ws.Range("J139").Value = 79994
$ws.Range("L139").Value = 79994
$ws.Range("N139").Value = -90274

$ws = $wb.Worksheets.Item("WVR")
# Row 54: 3413
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51040
# Row 113: 27752
$ws.Range("H113").Value = 976
$ws.Range("I113").Value = 798.2727
$ws.Range("J113").Value = 1106.3334
$ws.Range("K113").Value = 2394.8181
$ws.Range("L113").Value = 3319.0002
$ws.Range("M113").Value = -224.8181
$ws.Range("N113").Value = -7659.0002
